$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.095.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.818.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").Value = "'337.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "'0.9986"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'0.4302"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +12.62%  "
$ws.Range("D8").Value = "'0.3514"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.76%  "
$ws.Range("D9").Value = "'45.62"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").Value = "'1.152"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.81%  "
$ws.Range("D11").Value = "'0.07448"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.74%  "
$ws.Range("D12").Value = "'23.03"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("E13").Value = "  -0.15%  "
$ws.Range("D14").Value = "'6.265"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'7.281"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.83%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'1.813.71"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("E17").Value = "  +0.91%  "
$ws.Range("D18").Value = "'0.06687"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "'82.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "'6.486"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "'17.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "'28.127.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").Value = "'12.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").Value = "'2.396"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("E26").Value = "  +3.30%  "
$ws.Range("D27").Value = "'20.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").Value = "'156.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").Value = "'2.027.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.35%  "
$ws.Range("D30").Value = "'1.304"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -10.34%  "
$ws.Range("D31").Value = "'132.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("D32").Value = "'4.055"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").Value = "'5.973"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.55%  "
$ws.Range("D34").Value = "'0.09229"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.04%  "
$ws.Range("E35").Value = "  -2.76%  "
$ws.Range("D36").Value = "'0.02368"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("D37").Value = "'0.6744"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").Value = "'5.250"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("D39").Value = "'0.06259"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("D40").Value = "'0.2167"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("D41").Value = "'1.491"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").Value = "'1.219"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("D43").Value = "'8.217"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "'0.9987"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "'14.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("D46").Value = "'3.874"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Value = "'0.6137"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.10%  "
$ws.Range("D48").Value = "'128.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.31%  "
$ws.Range("D49").Value = "'2.048"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("E50").Value = "  -2.51%  "
$ws.Range("D51").Value = "'0.07113"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.90%  "
